$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two brand-new simulation rows (30, 31) need column-A styling (bold, thin border, centered)
# matching the existing index column; apply after setting values so cellXf matches exactly.

# Row 4: Holden
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Holden"
$ws.Cells.Item(4, 3).Value = 1.244954012271474
$ws.Cells.Item(4, 4).Value = 0.8351683804023381
$ws.Cells.Item(4, 5).Value = 0.9501010573252818
$ws.Cells.Item(4, 6).Value = 0.8052372103193526
$ws.Cells.Item(4, 7).Value = 0.7812762592089859
$ws.Cells.Item(4, 8).Value = 0.7812762592089859
$ws.Cells.Item(4, 9).Value = 0.7812762592089859
$ws.Cells.Item(4, 10).Value = 0.8337315578664718
$ws.Cells.Item(4, 11).Value = 4.091170572473485
$ws.Cells.Item(4, 12).Value = 4.091170572473485
$ws.Cells.Item(4, 13).Value = 0.948892488031977
$ws.Cells.Item(4, 14).Value = 0.7812762592089859
$ws.Cells.Item(4, 15).Value = 0.8337315578664718
$ws.Cells.Item(4, 16).Value = 2.462451065169978
$ws.Cells.Item(4, 17).Value = 0.8919163075958768
$ws.Cells.Item(4, 18).Value = 1.902059463182981
$ws.Cells.Item(4, 19).Value = 1.958334395888413
$ws.Cells.Item(4, 20).Value = 1.902059463182981
$ws.Cells.Item(4, 21).Value = 1.664069861718556
$ws.Cells.Item(4, 22).Value = 1.487511141216642
$ws.Cells.Item(4, 23).Value = 1.311316442237421

# Row 5: Rizzie Spiral
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$ws.Cells.Item(5, 3).Value = 0.2394496837162838
$ws.Cells.Item(5, 4).Value = 0.00112344790160153
$ws.Cells.Item(5, 5).Value = 0.002512355825932856
$ws.Cells.Item(5, 6).Value = 0.02752177195081789
$ws.Cells.Item(5, 7).Value = -0.001327213979300888
$ws.Cells.Item(5, 8).Value = -0.001327213979300888
$ws.Cells.Item(5, 9).Value = -0.001327213979300888
$ws.Cells.Item(5, 10).Value = 3.766703332693279
$ws.Cells.Item(5, 11).Value = 5.488704988933891
$ws.Cells.Item(5, 12).Value = 5.488704988933891
$ws.Cells.Item(5, 13).Value = 3.762758854977094
$ws.Cells.Item(5, 14).Value = -0.001327213979300888
$ws.Cells.Item(5, 15).Value = 3.766703332693279
$ws.Cells.Item(5, 16).Value = 4.627704160813585
$ws.Cells.Item(5, 17).Value = 1.884607844259606
$ws.Cells.Item(5, 18).Value = 3.08469370254929
$ws.Cells.Item(5, 19).Value = 3.085973559151034
$ws.Cells.Item(5, 20).Value = 3.08469370254929
$ws.Cells.Item(5, 21).Value = 2.314148365868451
$ws.Cells.Item(5, 22).Value = 1.851053249898901
$ws.Cells.Item(5, 23).Value = 1.66093090275245

# Row 6: RotRing OmegaMax-90
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "RotRing OmegaMax-90"
$ws.Cells.Item(6, 3).Value = 1.271305302961586
$ws.Cells.Item(6, 4).Value = 0.8157452884046038
$ws.Cells.Item(6, 5).Value = 0.9326485799463847
$ws.Cells.Item(6, 6).Value = 0.7935572696759775
$ws.Cells.Item(6, 7).Value = 0.7628010176520446
$ws.Cells.Item(6, 8).Value = 0.7628010176520446
$ws.Cells.Item(6, 9).Value = 0.7628010176520446
$ws.Cells.Item(6, 10).Value = 0.8780196593461093
$ws.Cells.Item(6, 11).Value = 3.133541325624553
$ws.Cells.Item(6, 12).Value = 3.133541325624553
$ws.Cells.Item(6, 13).Value = 0.9901082499467646
$ws.Cells.Item(6, 14).Value = 0.7628010176520446
$ws.Cells.Item(6, 15).Value = 0.8780196593461093
$ws.Cells.Item(6, 16).Value = 2.005780492485331
$ws.Cells.Item(6, 17).Value = 0.905334119646247
$ws.Cells.Item(6, 18).Value = 1.591454000874235
$ws.Cells.Item(6, 19).Value = 1.648069854972349
$ws.Cells.Item(6, 20).Value = 1.591454000874235
$ws.Cells.Item(6, 21).Value = 1.426752645642273
$ws.Cells.Item(6, 22).Value = 1.293962320044227
$ws.Cells.Item(6, 23).Value = 1.197215836694753

# Row 7: Equal Angle
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Equal Angle"
$ws.Cells.Item(7, 3).Value = 1.177227910082598
$ws.Cells.Item(7, 4).Value = 0.7814294200670343
$ws.Cells.Item(7, 5).Value = 0.8892397304736644
$ws.Cells.Item(7, 6).Value = 0.7546233449214304
$ws.Cells.Item(7, 7).Value = 0.7312364429001235
$ws.Cells.Item(7, 8).Value = 0.7312364429001235
$ws.Cells.Item(7, 9).Value = 0.7312364429001235
$ws.Cells.Item(7, 10).Value = 0.9022677565366983
$ws.Cells.Item(7, 11).Value = 5.697991800285292
$ws.Cells.Item(7, 12).Value = 5.697991800285292
$ws.Cells.Item(7, 13).Value = 1.010453874806119
$ws.Cells.Item(7, 14).Value = 0.7312364429001235
$ws.Cells.Item(7, 15).Value = 0.9022677565366983
$ws.Cells.Item(7, 16).Value = 3.300129778410995
$ws.Cells.Item(7, 17).Value = 0.8957537435051814
$ws.Cells.Item(7, 18).Value = 2.443831999907371
$ws.Cells.Item(7, 19).Value = 2.496499762431885
$ws.Cells.Item(7, 20).Value = 2.443831999907371
$ws.Cells.Item(7, 21).Value = 2.055183932548945
$ws.Cells.Item(7, 22).Value = 1.79039443461918
$ws.Cells.Item(7, 23).Value = 1.49305878500912

# Row 8: Tilt Rotate
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "Tilt Rotate"
$ws.Cells.Item(8, 3).Value = 0.6503803038767374
$ws.Cells.Item(8, 4).Value = 0.4143490544170742
$ws.Cells.Item(8, 5).Value = 0.476978684855654
$ws.Cells.Item(8, 6).Value = 0.4064132418113491
$ws.Cells.Item(8, 7).Value = 0.394274429213854
$ws.Cells.Item(8, 8).Value = 0.394274429213854
$ws.Cells.Item(8, 9).Value = 0.394274429213854
$ws.Cells.Item(8, 10).Value = 0.4086466486231076
$ws.Cells.Item(8, 11).Value = 32.38979835641413
$ws.Cells.Item(8, 12).Value = 32.38979835641413
$ws.Cells.Item(8, 13).Value = 0.4669330274124808
$ws.Cells.Item(8, 14).Value = 0.394274429213854
$ws.Cells.Item(8, 15).Value = 0.4086466486231076
$ws.Cells.Item(8, 16).Value = 16.39922250251862
$ws.Cells.Item(8, 17).Value = 0.4428126667393808
$ws.Cells.Item(8, 18).Value = 11.06423981141703
$ws.Cells.Item(8, 19).Value = 11.09180789663096
$ws.Cells.Item(8, 20).Value = 11.06423981141703
$ws.Cells.Item(8, 21).Value = 8.417424529776685
$ws.Cells.Item(8, 22).Value = 6.812794509664118
$ws.Cells.Item(8, 23).Value = 4.450971718328049

# Row 9: CLR
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "CLR"
$ws.Cells.Item(9, 3).Value = 0.9694832973219962
$ws.Cells.Item(9, 4).Value = 0.9496169627056327
$ws.Cells.Item(9, 5).Value = 0.9531403875924397
$ws.Cells.Item(9, 6).Value = 0.9550730015016526
$ws.Cells.Item(9, 7).Value = 0.9489874802718845
$ws.Cells.Item(9, 8).Value = 0.9489874802718845
$ws.Cells.Item(9, 9).Value = 0.9489874802718845
$ws.Cells.Item(9, 10).Value = 1.108154912647222
$ws.Cells.Item(9, 11).Value = 1.277531034035431
$ws.Cells.Item(9, 12).Value = 1.277531034035431
$ws.Cells.Item(9, 13).Value = 1.114145282339764
$ws.Cells.Item(9, 14).Value = 0.9489874802718845
$ws.Cells.Item(9, 15).Value = 1.108154912647222
$ws.Cells.Item(9, 16).Value = 1.192842973341327
$ws.Cells.Item(9, 17).Value = 1.030647650119831
$ws.Cells.Item(9, 18).Value = 1.111557808984846
$ws.Cells.Item(9, 19).Value = 1.112942111425031
$ws.Cells.Item(9, 20).Value = 1.111557808984846
$ws.Cells.Item(9, 21).Value = 1.071953453636744
$ws.Cells.Item(9, 22).Value = 1.047360258963772
$ws.Cells.Item(9, 23).Value = 1.034516544802003

# Row 10: Rizzie Hex
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Rizzie Hex"
$ws.Cells.Item(10, 3).Value = 1.000698421633131
$ws.Cells.Item(10, 4).Value = 0.993758846167905
$ws.Cells.Item(10, 5).Value = 0.9930648890227208
$ws.Cells.Item(10, 6).Value = 0.9909018666916163
$ws.Cells.Item(10, 7).Value = 0.9931101478651673
$ws.Cells.Item(10, 8).Value = 0.9931101478651673
$ws.Cells.Item(10, 9).Value = 0.9931101478651673
$ws.Cells.Item(10, 10).Value = 1.038860998658936
$ws.Cells.Item(10, 11).Value = 1.006463253683541
$ws.Cells.Item(10, 12).Value = 1.006463253683541
$ws.Cells.Item(10, 13).Value = 0.9879946921121984
$ws.Cells.Item(10, 14).Value = 0.9931101478651673
$ws.Cells.Item(10, 15).Value = 1.038860998658936
$ws.Cells.Item(10, 16).Value = 1.022662126171238
$ws.Cells.Item(10, 17).Value = 1.015962943840828
$ws.Cells.Item(10, 18).Value = 1.012811466735881
$ws.Cells.Item(10, 19).Value = 1.012796380455066
$ws.Cells.Item(10, 20).Value = 1.012811466735881
$ws.Cells.Item(10, 21).Value = 1.007874822307591
$ws.Cells.Item(10, 22).Value = 1.004921887419106
$ws.Cells.Item(10, 23).Value = 1.000606639479402

# Row 11: Matthies Hex
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Matthies Hex"
$ws.Cells.Item(11, 3).Value = 0.9518818160728856
$ws.Cells.Item(11, 4).Value = 0.9711807209918591
$ws.Cells.Item(11, 5).Value = 0.8586229122657004
$ws.Cells.Item(11, 6).Value = 0.8849646903713391
$ws.Cells.Item(11, 7).Value = 0.9685179208066926
$ws.Cells.Item(11, 8).Value = 0.9685179208066926
$ws.Cells.Item(11, 9).Value = 0.9685179208066926
$ws.Cells.Item(11, 10).Value = 1.205893345514888
$ws.Cells.Item(11, 11).Value = 1.671977490293228
$ws.Cells.Item(11, 12).Value = 1.671977490293228
$ws.Cells.Item(11, 13).Value = 1.229885850192874
$ws.Cells.Item(11, 14).Value = 0.9685179208066926
$ws.Cells.Item(11, 15).Value = 1.205893345514888
$ws.Cells.Item(11, 16).Value = 1.438935417904058
$ws.Cells.Item(11, 17).Value = 1.032258128890294
$ws.Cells.Item(11, 18).Value = 1.28212958553827
$ws.Cells.Item(11, 19).Value = 1.245497916024605
$ws.Cells.Item(11, 20).Value = 1.28212958553827
$ws.Cells.Item(11, 21).Value = 1.176252917220127
$ws.Cells.Item(11, 22).Value = 1.13470591793744
$ws.Cells.Item(11, 23).Value = 1.092865593313683

# Row 12: Tilt Rotate_Partial
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Tilt Rotate_Partial"
$ws.Cells.Item(12, 3).Value = 0.6365323209873008
$ws.Cells.Item(12, 4).Value = 0.3132099058101704
$ws.Cells.Item(12, 5).Value = 0.3814906670904826
$ws.Cells.Item(12, 6).Value = 0.4440747210660893
$ws.Cells.Item(12, 7).Value = 0.5149012431344437
$ws.Cells.Item(12, 8).Value = 0.5149012431344437
$ws.Cells.Item(12, 9).Value = 0.5149012431344437
$ws.Cells.Item(12, 10).Value = 0.4539281208841073
$ws.Cells.Item(12, 11).Value = 33.07672559460195
$ws.Cells.Item(12, 12).Value = 33.07672559460195
$ws.Cells.Item(12, 13).Value = 0.4060912466534862
$ws.Cells.Item(12, 14).Value = 0.5149012431344437
$ws.Cells.Item(12, 15).Value = 0.4539281208841073
$ws.Cells.Item(12, 16).Value = 16.76532685774303
$ws.Cells.Item(12, 17).Value = 0.417709393987295
$ws.Cells.Item(12, 18).Value = 11.34851831954017
$ws.Cells.Item(12, 19).Value = 11.30404812752551
$ws.Cells.Item(12, 20).Value = 11.34851831954017
$ws.Cells.Item(12, 21).Value = 8.606761406427747
$ws.Cells.Item(12, 22).Value = 6.988389373769087
$ws.Cells.Item(12, 23).Value = 4.528369227528504

# Row 13: RotRing OmegaMax-60
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "RotRing OmegaMax-60"
$ws.Cells.Item(13, 3).Value = 1.275742217209513
$ws.Cells.Item(13, 4).Value = 1.00183046425107
$ws.Cells.Item(13, 5).Value = 1.006496178566581
$ws.Cells.Item(13, 6).Value = 0.7934671790204072
$ws.Cells.Item(13, 7).Value = 1.126255008132724
$ws.Cells.Item(13, 8).Value = 1.126255008132724
$ws.Cells.Item(13, 9).Value = 1.126255008132724
$ws.Cells.Item(13, 10).Value = 0.7522912816772533
$ws.Cells.Item(13, 11).Value = 2.934682143652594
$ws.Cells.Item(13, 12).Value = 2.934682143652594
$ws.Cells.Item(13, 13).Value = 0.7582343430532114
$ws.Cells.Item(13, 14).Value = 1.126255008132724
$ws.Cells.Item(13, 15).Value = 0.7522912816772533
$ws.Cells.Item(13, 16).Value = 1.843486712664923
$ws.Cells.Item(13, 17).Value = 0.8793937301219171
$ws.Cells.Item(13, 18).Value = 1.604409477820857
$ws.Cells.Item(13, 19).Value = 1.564489867965476
$ws.Cells.Item(13, 20).Value = 1.604409477820857
$ws.Cells.Item(13, 21).Value = 1.454931153007288
$ws.Cells.Item(13, 22).Value = 1.389195924032375
$ws.Cells.Item(13, 23).Value = 1.206124851945419

# Row 14: Equal Angle_Partial
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "Equal Angle_Partial"
$ws.Cells.Item(14, 3).Value = 1.117018488312304
$ws.Cells.Item(14, 4).Value = 0.6336167774653726
$ws.Cells.Item(14, 5).Value = 0.694107369745399
$ws.Cells.Item(14, 6).Value = 0.8032576241514495
$ws.Cells.Item(14, 7).Value = 1.042194791570058
$ws.Cells.Item(14, 8).Value = 1.042194791570058
$ws.Cells.Item(14, 9).Value = 1.042194791570058
$ws.Cells.Item(14, 10).Value = 0.8917008121109385
$ws.Cells.Item(14, 11).Value = 7.371875189357552
$ws.Cells.Item(14, 12).Value = 7.371875189357552
$ws.Cells.Item(14, 13).Value = 0.8396171408691885
$ws.Cells.Item(14, 14).Value = 1.042194791570058
$ws.Cells.Item(14, 15).Value = 0.8917008121109385
$ws.Cells.Item(14, 16).Value = 4.131788000734245
$ws.Cells.Item(14, 17).Value = 0.7929040909281688
$ws.Cells.Item(14, 18).Value = 3.101923597679516
$ws.Cells.Item(14, 19).Value = 2.985894457071296
$ws.Cells.Item(14, 20).Value = 3.101923597679516
$ws.Cells.Item(14, 21).Value = 2.499969540695987
$ws.Cells.Item(14, 22).Value = 2.208414590870801
$ws.Cells.Item(14, 23).Value = 1.674173524197783

# Row 15: Rizzie Hex_Partial
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rizzie Hex_Partial"
$ws.Cells.Item(15, 3).Value = 0.6672487745394171
$ws.Cells.Item(15, 4).Value = 0.8056482910811112
$ws.Cells.Item(15, 5).Value = 0.6635741383564324
$ws.Cells.Item(15, 6).Value = 1.320002956595297
$ws.Cells.Item(15, 7).Value = 1.67386080182847
$ws.Cells.Item(15, 8).Value = 1.67386080182847
$ws.Cells.Item(15, 9).Value = 1.67386080182847
$ws.Cells.Item(15, 10).Value = 1.358845307334478
$ws.Cells.Item(15, 11).Value = 0.7361896534508522
$ws.Cells.Item(15, 12).Value = 0.7361896534508522
$ws.Cells.Item(15, 13).Value = 0.8331818604903317
$ws.Cells.Item(15, 14).Value = 1.67386080182847
$ws.Cells.Item(15, 15).Value = 1.358845307334478
$ws.Cells.Item(15, 16).Value = 1.047517480392665
$ws.Cells.Item(15, 17).Value = 1.011209722845455
$ws.Cells.Item(15, 18).Value = 1.256298587537933
$ws.Cells.Item(15, 19).Value = 0.9195363663805874
$ws.Cells.Item(15, 20).Value = 1.256298587537933
$ws.Cells.Item(15, 21).Value = 1.108117475242558
$ws.Cells.Item(15, 22).Value = 1.22126614055974
$ws.Cells.Item(15, 23).Value = 1.007318972959549

# Row 16: ND Single
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "ND Single"
$ws.Cells.Item(16, 3).Value = 0.04955152399999998
$ws.Cells.Item(16, 4).Value = -0.006433271399999997
$ws.Cells.Item(16, 5).Value = 0.002693453900000002
$ws.Cells.Item(16, 6).Value = 0.007878550500000006
$ws.Cells.Item(16, 7).Value = 0.006236068099999995
$ws.Cells.Item(16, 8).Value = 0.006236068099999995
$ws.Cells.Item(16, 9).Value = 0.006236068099999995
$ws.Cells.Item(16, 10).Value = 0.003096941800000001
$ws.Cells.Item(16, 11).Value = 60.75982100000006
$ws.Cells.Item(16, 12).Value = 60.75982100000006
$ws.Cells.Item(16, 13).Value = 0.003041090100000003
$ws.Cells.Item(16, 14).Value = 0.006236068099999995
$ws.Cells.Item(16, 15).Value = 0.003096941800000001
$ws.Cells.Item(16, 16).Value = 30.38145897090003
$ws.Cells.Item(16, 17).Value = 0.002895197850000002
$ws.Cells.Item(16, 18).Value = 20.25638466996669
$ws.Cells.Item(16, 19).Value = 20.25520379856669
$ws.Cells.Item(16, 20).Value = 20.25638466996669
$ws.Cells.Item(16, 21).Value = 15.19296186595002
$ws.Cells.Item(16, 22).Value = 12.15561670638001
$ws.Cells.Item(16, 23).Value = 7.603235669625007

# Row 17: RD Single
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "RD Single"
$ws.Cells.Item(17, 3).Value = 0.049551524
$ws.Cells.Item(17, 4).Value = -0.0064332714
$ws.Cells.Item(17, 5).Value = 0.0026934539
$ws.Cells.Item(17, 6).Value = 0.0078785505
$ws.Cells.Item(17, 7).Value = 0.0062360681
$ws.Cells.Item(17, 8).Value = 0.0062360681
$ws.Cells.Item(17, 9).Value = 0.0062360681
$ws.Cells.Item(17, 10).Value = 0.0030969418
$ws.Cells.Item(17, 11).Value = 60.759821
$ws.Cells.Item(17, 12).Value = 60.759821
$ws.Cells.Item(17, 13).Value = 0.003041090099999999
$ws.Cells.Item(17, 14).Value = 0.0062360681
$ws.Cells.Item(17, 15).Value = 0.0030969418
$ws.Cells.Item(17, 16).Value = 30.3814589709
$ws.Cells.Item(17, 17).Value = 0.00289519785
$ws.Cells.Item(17, 18).Value = 20.25638466996667
$ws.Cells.Item(17, 19).Value = 20.25520379856667
$ws.Cells.Item(17, 20).Value = 20.25638466996667
$ws.Cells.Item(17, 21).Value = 15.19296186595
$ws.Cells.Item(17, 22).Value = 12.15561670638
$ws.Cells.Item(17, 23).Value = 7.603235669625

# Row 18: TD Single
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "TD Single"
$ws.Cells.Item(18, 3).Value = 0.049551524
$ws.Cells.Item(18, 4).Value = -0.0064332714
$ws.Cells.Item(18, 5).Value = 0.0026934539
$ws.Cells.Item(18, 6).Value = 0.0078785505
$ws.Cells.Item(18, 7).Value = 0.0062360681
$ws.Cells.Item(18, 8).Value = 0.0062360681
$ws.Cells.Item(18, 9).Value = 0.0062360681
$ws.Cells.Item(18, 10).Value = 0.0030969418
$ws.Cells.Item(18, 11).Value = 60.759821
$ws.Cells.Item(18, 12).Value = 60.759821
$ws.Cells.Item(18, 13).Value = 0.003041090099999999
$ws.Cells.Item(18, 14).Value = 0.0062360681
$ws.Cells.Item(18, 15).Value = 0.0030969418
$ws.Cells.Item(18, 16).Value = 30.3814589709
$ws.Cells.Item(18, 17).Value = 0.00289519785
$ws.Cells.Item(18, 18).Value = 20.25638466996667
$ws.Cells.Item(18, 19).Value = 20.25520379856667
$ws.Cells.Item(18, 20).Value = 20.25638466996667
$ws.Cells.Item(18, 21).Value = 15.19296186595
$ws.Cells.Item(18, 22).Value = 12.15561670638
$ws.Cells.Item(18, 23).Value = 7.603235669625

# Row 19: Morris Single
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "Morris Single"
$ws.Cells.Item(19, 3).Value = -0.0082885554
$ws.Cells.Item(19, 4).Value = 0.0035140694
$ws.Cells.Item(19, 5).Value = -0.0063507761
$ws.Cells.Item(19, 6).Value = 0.0016794768
$ws.Cells.Item(19, 7).Value = 0.0002567862
$ws.Cells.Item(19, 8).Value = 0.0002567862
$ws.Cells.Item(19, 9).Value = 0.0002567862
$ws.Cells.Item(19, 10).Value = 0.056027822
$ws.Cells.Item(19, 11).Value = 0.0013547887
$ws.Cells.Item(19, 12).Value = 0.0013547887
$ws.Cells.Item(19, 13).Value = 10.971864
$ws.Cells.Item(19, 14).Value = 0.0002567862
$ws.Cells.Item(19, 15).Value = 0.056027822
$ws.Cells.Item(19, 16).Value = 0.02869130535
$ws.Cells.Item(19, 17).Value = 0.02483852295
$ws.Cells.Item(19, 18).Value = 0.0192131323
$ws.Cells.Item(19, 19).Value = 0.01701061153333333
$ws.Cells.Item(19, 20).Value = 0.0192131323
$ws.Cells.Item(19, 21).Value = 0.0128221552
$ws.Cells.Item(19, 22).Value = 0.0103090814
$ws.Cells.Item(19, 23).Value = 1.37750720145

# Row 20: Ring Perpendicular to ND
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "Ring Perpendicular to ND"
$ws.Cells.Item(20, 3).Value = 0.2358865764986303
$ws.Cells.Item(20, 4).Value = 0.000972024362191781
$ws.Cells.Item(20, 5).Value = 0.002328569554794522
$ws.Cells.Item(20, 6).Value = 0.0261257860890411
$ws.Cells.Item(20, 7).Value = -0.001174491636986302
$ws.Cells.Item(20, 8).Value = -0.001174491636986302
$ws.Cells.Item(20, 9).Value = -0.001174491636986302
$ws.Cells.Item(20, 10).Value = 3.546992953636714
$ws.Cells.Item(20, 11).Value = 7.927111028986301
$ws.Cells.Item(20, 12).Value = 7.927111028986301
$ws.Cells.Item(20, 13).Value = 3.545767832329591
$ws.Cells.Item(20, 14).Value = -0.001174491636986302
$ws.Cells.Item(20, 15).Value = 3.546992953636714
$ws.Cells.Item(20, 16).Value = 5.737051991311508
$ws.Cells.Item(20, 17).Value = 1.774660761595754
$ws.Cells.Item(20, 18).Value = 3.824309830328676
$ws.Cells.Item(20, 19).Value = 3.825477517392603
$ws.Cells.Item(20, 20).Value = 3.824309830328676
$ws.Cells.Item(20, 21).Value = 2.868814515135206
$ws.Cells.Item(20, 22).Value = 2.294816713780768
$ws.Cells.Item(20, 23).Value = 1.910501284977535

# Row 21: Ring Perpendicular to RD
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "Ring Perpendicular to RD"
$ws.Cells.Item(21, 3).Value = 0.2285312454789473
$ws.Cells.Item(21, 4).Value = 0.0006797100494736845
$ws.Cells.Item(21, 5).Value = 0.002342972905263157
$ws.Cells.Item(21, 6).Value = 0.02540550047368421
$ws.Cells.Item(21, 7).Value = -0.0008819695421052635
$ws.Cells.Item(21, 8).Value = -0.0008819695421052635
$ws.Cells.Item(21, 9).Value = -0.0008819695421052635
$ws.Cells.Item(21, 10).Value = 3.407102321590526
$ws.Cells.Item(21, 11).Value = 10.01261273836842
$ws.Cells.Item(21, 12).Value = 10.01261273836842
$ws.Cells.Item(21, 13).Value = 3.405923355662632
$ws.Cells.Item(21, 14).Value = -0.0008819695421052635
$ws.Cells.Item(21, 15).Value = 3.407102321590526
$ws.Cells.Item(21, 16).Value = 6.709857529979473
$ws.Cells.Item(21, 17).Value = 1.704722647247894
$ws.Cells.Item(21, 18).Value = 4.472944363472281
$ws.Cells.Item(21, 19).Value = 4.474019344288069
$ws.Cells.Item(21, 20).Value = 4.47294436347228
$ws.Cells.Item(21, 21).Value = 3.355294015830526
$ws.Cells.Item(21, 22).Value = 2.684058818755999
$ws.Cells.Item(21, 23).Value = 2.135214484373355

# Row 22: Ring Perpendicular to TD
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "Ring Perpendicular to TD"
$ws.Cells.Item(22, 3).Value = 0.2285312454789473
$ws.Cells.Item(22, 4).Value = 0.0006797100494736845
$ws.Cells.Item(22, 5).Value = 0.002342972905263157
$ws.Cells.Item(22, 6).Value = 0.02540550047368421
$ws.Cells.Item(22, 7).Value = -0.0008819695421052635
$ws.Cells.Item(22, 8).Value = -0.0008819695421052635
$ws.Cells.Item(22, 9).Value = -0.0008819695421052635
$ws.Cells.Item(22, 10).Value = 3.407102321590526
$ws.Cells.Item(22, 11).Value = 10.01261273836842
$ws.Cells.Item(22, 12).Value = 10.01261273836842
$ws.Cells.Item(22, 13).Value = 3.405923355662632
$ws.Cells.Item(22, 14).Value = -0.0008819695421052635
$ws.Cells.Item(22, 15).Value = 3.407102321590526
$ws.Cells.Item(22, 16).Value = 6.709857529979473
$ws.Cells.Item(22, 17).Value = 1.704722647247894
$ws.Cells.Item(22, 18).Value = 4.472944363472281
$ws.Cells.Item(22, 19).Value = 4.474019344288069
$ws.Cells.Item(22, 20).Value = 4.47294436347228
$ws.Cells.Item(22, 21).Value = 3.355294015830526
$ws.Cells.Item(22, 22).Value = 2.684058818755999
$ws.Cells.Item(22, 23).Value = 2.135214484373355

# Row 23: OffsetFTD
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "OffsetFTD"
$ws.Cells.Item(23, 3).Value = 0.0287933131695283
$ws.Cells.Item(23, 4).Value = 0.03001048724434661
$ws.Cells.Item(23, 5).Value = 0.08375329453524756
$ws.Cells.Item(23, 6).Value = 2.215396312838418
$ws.Cells.Item(23, 7).Value = 0.0003904390819234697
$ws.Cells.Item(23, 8).Value = 0.0003904390819234697
$ws.Cells.Item(23, 9).Value = 0.0003904390819234697
$ws.Cells.Item(23, 10).Value = 3.835382676343219
$ws.Cells.Item(23, 11).Value = 0.001195947500271863
$ws.Cells.Item(23, 12).Value = 0.001195947500271863
$ws.Cells.Item(23, 13).Value = 1.719955642120351
$ws.Cells.Item(23, 14).Value = 0.0003904390819234697
$ws.Cells.Item(23, 15).Value = 3.835382676343219
$ws.Cells.Item(23, 16).Value = 1.918289311921745
$ws.Cells.Item(23, 17).Value = 1.959567985439233
$ws.Cells.Item(23, 18).Value = 1.278989687641805
$ws.Cells.Item(23, 19).Value = 1.306777306126246
$ws.Cells.Item(23, 20).Value = 1.278989687641805
$ws.Cells.Item(23, 21).Value = 0.9801805893651654
$ws.Cells.Item(23, 22).Value = 0.784222559308517
$ws.Cells.Item(23, 23).Value = 0.989359764104163

# Row 24: OffsetATD
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "OffsetATD"
$ws.Cells.Item(24, 3).Value = 1.978017528762467
$ws.Cells.Item(24, 4).Value = 2.103189352075622
$ws.Cells.Item(24, 5).Value = 1.930752885397657
$ws.Cells.Item(24, 6).Value = 0.3262007917718728
$ws.Cells.Item(24, 7).Value = 0.3335060775395657
$ws.Cells.Item(24, 8).Value = 0.3335060775395657
$ws.Cells.Item(24, 9).Value = 0.3335060775395657
$ws.Cells.Item(24, 10).Value = 0.0062426259214741
$ws.Cells.Item(24, 11).Value = 0.4867331853026556
$ws.Cells.Item(24, 12).Value = 0.4867331853026556
$ws.Cells.Item(24, 13).Value = 0.7639804803541627
$ws.Cells.Item(24, 14).Value = 0.3335060775395657
$ws.Cells.Item(24, 15).Value = 0.0062426259214741
$ws.Cells.Item(24, 16).Value = 0.2464879056120649
$ws.Cells.Item(24, 17).Value = 0.9684977556595658
$ws.Cells.Item(24, 18).Value = 0.2754939629212318
$ws.Cells.Item(24, 19).Value = 0.8079095655405957
$ws.Cells.Item(24, 20).Value = 0.2754939629212318
$ws.Cells.Item(24, 21).Value = 0.6893086935403382
$ws.Cells.Item(24, 22).Value = 0.6181481703401837
$ws.Cells.Item(24, 23).Value = 0.9910778658906845

# Row 25: OffsetF45
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "OffsetF45"
$ws.Cells.Item(25, 3).Value = 0.569984070182464
$ws.Cells.Item(25, 4).Value = 0.3632148420855516
$ws.Cells.Item(25, 5).Value = 1.916597466133352
$ws.Cells.Item(25, 6).Value = 0.5588256154470996
$ws.Cells.Item(25, 7).Value = 1.379937255016048
$ws.Cells.Item(25, 8).Value = 1.379937255016048
$ws.Cells.Item(25, 9).Value = 1.379937255016048
$ws.Cells.Item(25, 10).Value = 0.002237596383924254
$ws.Cells.Item(25, 11).Value = 0.001328490850010016
$ws.Cells.Item(25, 12).Value = 0.001328490850010016
$ws.Cells.Item(25, 13).Value = 0.08313614211869384
$ws.Cells.Item(25, 14).Value = 1.379937255016048
$ws.Cells.Item(25, 15).Value = 0.002237596383924254
$ws.Cells.Item(25, 16).Value = 0.001783043616967135
$ws.Cells.Item(25, 17).Value = 0.9594175312586382
$ws.Cells.Item(25, 18).Value = 0.4611677807499943
$ws.Cells.Item(25, 19).Value = 0.6400545177890954
$ws.Cells.Item(25, 20).Value = 0.4611677807499943
$ws.Cells.Item(25, 21).Value = 0.8250252020958337
$ws.Cells.Item(25, 22).Value = 0.9360076126798766
$ws.Cells.Item(25, 23).Value = 0.609407684777143

# Row 26: OffsetA45
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "OffsetA45"
$ws.Cells.Item(26, 3).Value = 0.3930322337325861
$ws.Cells.Item(26, 4).Value = 0.4028156866737078
$ws.Cells.Item(26, 5).Value = 0.005356301082907026
$ws.Cells.Item(26, 6).Value = 1.128610781007443
$ws.Cells.Item(26, 7).Value = 1.376818959164313
$ws.Cells.Item(26, 8).Value = 1.376818959164313
$ws.Cells.Item(26, 9).Value = 1.376818959164313
$ws.Cells.Item(26, 10).Value = 2.359611698722838
$ws.Cells.Item(26, 11).Value = 0.4869577846942586
$ws.Cells.Item(26, 12).Value = 0.4869577846942586
$ws.Cells.Item(26, 13).Value = 1.911625012920187
$ws.Cells.Item(26, 14).Value = 1.376818959164313
$ws.Cells.Item(26, 15).Value = 2.359611698722838
$ws.Cells.Item(26, 16).Value = 1.423284741708549
$ws.Cells.Item(26, 17).Value = 1.182483999902873
$ws.Cells.Item(26, 18).Value = 1.407796147527137
$ws.Cells.Item(26, 19).Value = 0.950641928166668
$ws.Cells.Item(26, 20).Value = 1.407796147527137
$ws.Cells.Item(26, 21).Value = 1.057186185916079
$ws.Cells.Item(26, 22).Value = 1.121112740565726
$ws.Cells.Item(26, 23).Value = 1.00810355724978

# Row 27: OffsetFRD
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "OffsetFRD"
$ws.Cells.Item(27, 3).Value = 0.02879331316952831
$ws.Cells.Item(27, 4).Value = 0.0300104872443466
$ws.Cells.Item(27, 5).Value = 0.08375329453524757
$ws.Cells.Item(27, 6).Value = 2.215396312838417
$ws.Cells.Item(27, 7).Value = 0.0003904390819234731
$ws.Cells.Item(27, 8).Value = 0.0003904390819234731
$ws.Cells.Item(27, 9).Value = 0.0003904390819234731
$ws.Cells.Item(27, 10).Value = 3.835382676343219
$ws.Cells.Item(27, 11).Value = 0.001195947486230521
$ws.Cells.Item(27, 12).Value = 0.001195947486230521
$ws.Cells.Item(27, 13).Value = 1.71995564212035
$ws.Cells.Item(27, 14).Value = 0.0003904390819234731
$ws.Cells.Item(27, 15).Value = 3.835382676343219
$ws.Cells.Item(27, 16).Value = 1.918289311914725
$ws.Cells.Item(27, 17).Value = 1.959567985439233
$ws.Cells.Item(27, 18).Value = 1.278989687637125
$ws.Cells.Item(27, 19).Value = 1.306777306121566
$ws.Cells.Item(27, 20).Value = 1.278989687637125
$ws.Cells.Item(27, 21).Value = 0.9801805893616552
$ws.Cells.Item(27, 22).Value = 0.7842225593057088
$ws.Cells.Item(27, 23).Value = 0.9893597641024079

# Row 28: OffsetARD
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "OffsetARD"
$ws.Cells.Item(28, 3).Value = 1.978017528762467
$ws.Cells.Item(28, 4).Value = 2.103189352075622
$ws.Cells.Item(28, 5).Value = 1.930752885397657
$ws.Cells.Item(28, 6).Value = 0.3262007917718726
$ws.Cells.Item(28, 7).Value = 0.3335060775395662
$ws.Cells.Item(28, 8).Value = 0.3335060775395662
$ws.Cells.Item(28, 9).Value = 0.3335060775395662
$ws.Cells.Item(28, 10).Value = 0.006242625921474081
$ws.Cells.Item(28, 11).Value = 0.4867331853026644
$ws.Cells.Item(28, 12).Value = 0.4867331853026644
$ws.Cells.Item(28, 13).Value = 0.7639804803541627
$ws.Cells.Item(28, 14).Value = 0.3335060775395662
$ws.Cells.Item(28, 15).Value = 0.006242625921474081
$ws.Cells.Item(28, 16).Value = 0.2464879056120692
$ws.Cells.Item(28, 17).Value = 0.9684977556595658
$ws.Cells.Item(28, 18).Value = 0.2754939629212349
$ws.Cells.Item(28, 19).Value = 0.8079095655405987
$ws.Cells.Item(28, 20).Value = 0.2754939629212349
$ws.Cells.Item(28, 21).Value = 0.6893086935403405
$ws.Cells.Item(28, 22).Value = 0.6181481703401857
$ws.Cells.Item(28, 23).Value = 0.9910778658906857

# Row 29: Gaussian Quadrature
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "Gaussian Quadrature"
$ws.Cells.Item(29, 3).Value = 0.04600138865378531
$ws.Cells.Item(29, 4).Value = 1.875846477088833
$ws.Cells.Item(29, 5).Value = 1.729017613154956
$ws.Cells.Item(29, 6).Value = 1.072011405849706
$ws.Cells.Item(29, 7).Value = 0.01508196230821052
$ws.Cells.Item(29, 8).Value = 0.01508196230821052
$ws.Cells.Item(29, 9).Value = 0.01508196230821052
$ws.Cells.Item(29, 10).Value = 0.06394563458397068
$ws.Cells.Item(29, 11).Value = 3.575753599251714
$ws.Cells.Item(29, 12).Value = 3.575753599251714
$ws.Cells.Item(29, 13).Value = 3.406828348519421
$ws.Cells.Item(29, 14).Value = 0.01508196230821052
$ws.Cells.Item(29, 15).Value = 0.06394563458397068
$ws.Cells.Item(29, 16).Value = 1.819849616917842
$ws.Cells.Item(29, 17).Value = 0.8964816238694632
$ws.Cells.Item(29, 18).Value = 1.218260398714632
$ws.Cells.Item(29, 19).Value = 1.789572282330213
$ws.Cells.Item(29, 20).Value = 1.218260398714632
$ws.Cells.Item(29, 21).Value = 1.345949702324713
$ws.Cells.Item(29, 22).Value = 1.079776154321412
$ws.Cells.Item(29, 23).Value = 1.473060803676325

# Row 30: Michael-CCHex
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "Michael-CCHex"
$ws.Cells.Item(30, 3).Value = 0.3884201130979834
$ws.Cells.Item(30, 4).Value = 1.152465689530825
$ws.Cells.Item(30, 5).Value = 0.2453017638869989
$ws.Cells.Item(30, 6).Value = 0.1365809871778847
$ws.Cells.Item(30, 7).Value = 0.08535652632257268
$ws.Cells.Item(30, 8).Value = 0.08535652632257268
$ws.Cells.Item(30, 9).Value = 0.08535652632257268
$ws.Cells.Item(30, 10).Value = 2.080813916858228
$ws.Cells.Item(30, 11).Value = 0.2709192540440629
$ws.Cells.Item(30, 12).Value = 0.2709192540440629
$ws.Cells.Item(30, 13).Value = 3.126665604690211
$ws.Cells.Item(30, 14).Value = 0.08535652632257268
$ws.Cells.Item(30, 15).Value = 2.080813916858228
$ws.Cells.Item(30, 16).Value = 1.175866585451146
$ws.Cells.Item(30, 17).Value = 1.163057840372614
$ws.Cells.Item(30, 18).Value = 0.812363232408288
$ws.Cells.Item(30, 19).Value = 0.86567831159643
$ws.Cells.Item(30, 20).Value = 0.812363232408288
$ws.Cells.Item(30, 21).Value = 0.6705978652779657
$ws.Cells.Item(30, 22).Value = 0.5535495974868871
$ws.Cells.Item(30, 23).Value = 0.9358154819510958

# Row 31: Michael-SNHex
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "Michael-SNHex"
$ws.Cells.Item(31, 3).Value = 0.003498197981038845
$ws.Cells.Item(31, 4).Value = 0.393863496498899
$ws.Cells.Item(31, 5).Value = 0.502440844649569
$ws.Cells.Item(31, 6).Value = 0.1306662516702358
$ws.Cells.Item(31, 7).Value = 0.1650157339284239
$ws.Cells.Item(31, 8).Value = 0.1650157339284239
$ws.Cells.Item(31, 9).Value = 0.1650157339284239
$ws.Cells.Item(31, 10).Value = 1.23444317543616
$ws.Cells.Item(31, 11).Value = -0.0001446655140792162
$ws.Cells.Item(31, 12).Value = -0.0001446655140792162
$ws.Cells.Item(31, 13).Value = 3.015894795148841
$ws.Cells.Item(31, 14).Value = 0.1650157339284239
$ws.Cells.Item(31, 15).Value = 1.23444317543616
$ws.Cells.Item(31, 16).Value = 0.6171492549610403
$ws.Cells.Item(31, 17).Value = 0.8684420100428645
$ws.Cells.Item(31, 18).Value = 0.4664380812835016
$ws.Cells.Item(31, 19).Value = 0.5789131181905499
$ws.Cells.Item(31, 20).Value = 0.4664380812835016
$ws.Cells.Item(31, 21).Value = 0.4754387721250184
$ws.Cells.Item(31, 22).Value = 0.4133541644856994
$ws.Cells.Item(31, 23).Value = 0.680709728724886

# Ensure the new rows (30, 31) carry the same index-column style as the rest of column A
for ($r = 30; $r -le 31; $r++) {
  $cell = $ws.Cells.Item($r, 1)
  $cell.Font.Bold = $true
  $cell.HorizontalAlignment = -4108
  $cell.VerticalAlignment = -4160
  $cell.Borders.LineStyle = 1
}
